$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Get-ParaIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    throw "Paragraph not found: $text"
}

function Get-ParaByText($text) {
    $idx = Get-ParaIndexByText $text
    return $d.Paragraphs($idx)
}

# Reliable whole-paragraph text replacement. Paragraph.Range.Text = "..."
# only touches the paragraph's first run when the paragraph holds more than
# one run (or a non-text element such as a symbol/page-break marker), so we
# explicitly delete the paragraph's content (excluding its end-of-paragraph
# mark) and insert fresh text instead.
function Set-ParaText($p, $text) {
    $r = $p.Range
    $body = $d.Range($r.Start, $r.End - 1)
    if ($body.End -gt $body.Start) {
        $body.Delete() | Out-Null
    }
    $ins = $d.Range($body.Start, $body.Start)
    $ins.InsertAfter($text)
}

# Inserts a new paragraph right after $p (inheriting $p's paragraph
# formatting/list level) containing $text, and returns the new Paragraph.
function Add-ParaAfter($p, $text) {
    $beforeText = $p.Range.Text.TrimEnd([char]13)
    $p.Range.InsertParagraphAfter()
    $idx = (Get-ParaIndexByText $beforeText) + 1
    $newPara = $d.Paragraphs($idx)
    Set-ParaText $newPara $text
    return $newPara
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-Text "Functions currently not working" "State of the S.C.O.R.E. Project"

# ---------------------------------------------------------------------------
# 2. Insert "Help" bullet after "Lecturer UI" (inherits ilvl 1 / numId 1)
# ---------------------------------------------------------------------------
Add-ParaAfter (Get-ParaByText "Lecturer UI") "Help" | Out-Null

# ---------------------------------------------------------------------------
# 3. "Algorithm – Currently in testing stages" gains a trailing clause
# ---------------------------------------------------------------------------
$pAlgoTesting = Get-ParaByText "Algorithm – Currently in testing stages"
$r = $pAlgoTesting.Range
$bodyEnd = $d.Range($r.Start, $r.End - 1).End
$insAfter = $d.Range($bodyEnd, $bodyEnd)
$insAfter.InsertAfter("; not yet functional")

# ---------------------------------------------------------------------------
# 4. "No functionality for the following users:" -> "... following roles:"
# ---------------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute(" following users:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $found.Find.Found) { throw "couldn't find 'following users:' text" }
$found.Text = " following roles:"

# ---------------------------------------------------------------------------
# 5. "Menu navigation bar does not work" -> split text, drop sub-bullets,
#    add "Help is not yet online" bullet
# ---------------------------------------------------------------------------
$pMenu = Get-ParaByText "Menu navigation bar does not work"
Set-ParaText $pMenu "Menu navigation bar works, but none of its targets are implemented"

# Delete the four now-obsolete sub-bullets: Course Chooser, Review Courses,
# Personal, Settings
$pCourseChooser = Get-ParaByText "Course Chooser"
$pSettings = Get-ParaByText "Settings"
$delStart = $pCourseChooser.Range.Start
$delEnd = $pSettings.Range.End
$d.Range($delStart, $delEnd).Delete() | Out-Null

# Add the replacement bullet after "Menu navigation bar ..." at the same
# list level (ilvl 1 / numId 2)
$pMenu = Get-ParaByText "Menu navigation bar works, but none of its targets are implemented"
Add-ParaAfter $pMenu "Help is not yet online" | Out-Null

# ---------------------------------------------------------------------------
# 6. UI -> Django -> Algorithm paragraph rewritten as plain sentence
# ---------------------------------------------------------------------------
$pUiDjango = Get-ParaByText "UIDjangoAlgorithm"
Set-ParaText $pUiDjango "The algorithm is not yet callable from Django; a stub is currently being used."

# ---------------------------------------------------------------------------
# 7. "Program Administrator can populate database ..." rewritten, and a new
#    bullet added below it.
# ---------------------------------------------------------------------------
$pProgAdmin = Get-ParaByText "Program Administrator can populate database with current user interface"
Set-ParaText $pProgAdmin "The authentication system works properly – users can log in."
$pProgAdmin = Get-ParaByText "The authentication system works properly – users can log in."
Add-ParaAfter $pProgAdmin "Program Administrator UI is in a state that conforms to the requirements document." | Out-Null

# ---------------------------------------------------------------------------
# 8. "Runs algorithm " paragraph split: trailing space becomes a new
#    "Adds users " bullet at the same list level.
# ---------------------------------------------------------------------------
$pRunsAlgo = Get-ParaByText "Runs algorithm "
Set-ParaText $pRunsAlgo "Runs algorithm"
$pRunsAlgo = Get-ParaByText "Runs algorithm"
Add-ParaAfter $pRunsAlgo "Adds users " | Out-Null

# ---------------------------------------------------------------------------
# 9. Database bullet rewritten
# ---------------------------------------------------------------------------
$pDatabase = Get-ParaByText "Database – All components of our database works and is able to accept values to populate the tables"
Set-ParaText $pDatabase "Database – All components of our database work, accept data, and return data. The database conforms to the specification – particularly the database design document."

# ---------------------------------------------------------------------------
# 10. Back-end bullet rewritten, plus a trailing empty ListParagraph
# ---------------------------------------------------------------------------
$pBackend = Get-ParaByText "Back-end – All files and programs are functionally working from remote server"
Set-ParaText $pBackend "Back-end – Lighttpd runs on our server, and calls Django properly. Django functions properly, and pages are served over the internet as expected."

# New trailing empty ListParagraph-styled paragraph at the very end
$pBackend = Get-ParaByText "Back-end – Lighttpd runs on our server, and calls Django properly. Django functions properly, and pages are served over the internet as expected."
$pBackend.Range.InsertParagraphAfter()
$trailingIdx = $d.Paragraphs.Count
$pTrailing = $d.Paragraphs($trailingIdx)
$pTrailing.Range.ListFormat.RemoveNumbers()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
